$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" -------------
# Overview sheet: columns E (zh-cn) and F (de-de), rows 2-4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# zh-cn sheet: Status column (C), rows 2-4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

# de-de sheet: Status column (C), rows 2-4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Shrink the now-shorter "In Translation" status columns -----------------
# The column previously sized for "Ready for handoff" (17 chars) is narrowed
# to fit "In Translation" (14 chars). Use the ColumnWidth nearest the target.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
